$d = $word.ActiveDocument

# The document currently has this bookmark sitting right after "repositories."
# at the end of the third paragraph. We need to relocate it to the very start
# of what is currently the second paragraph ("I wonder the difference..."),
# and then delete the whole first paragraph ("Webteam meetings are fun.").

# Remove the existing "_GoBack" bookmark (it will be re-added at its new spot).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create the bookmark as a zero-length mark at the start of the paragraph
# that currently follows the one we are about to delete. Doing this before the
# deletion (while that paragraph is not the very first paragraph in the
# document) keeps the bookmark cleanly inside a single paragraph instead of
# spanning a paragraph boundary.
$targetPara = $d.Paragraphs(2)
$bookmarkRange = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Delete the first paragraph (including its paragraph mark) entirely.
$d.Paragraphs(1).Range.Delete()
